$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values are written with a leading apostrophe to force
# Excel to store them as text, matching the original inlineStr
# (text) cell type for numeric-looking strings such as "22.500.32"
# or "20.90" (which would otherwise be auto-converted to numbers).
$ws.Range("D2").Value = "'22.500.32"
$ws.Range("E2").Value = "'  +0.08%  "
$ws.Range("D3").Value = "'1.573.17"
$ws.Range("E3").Value = "'  +0.01%  "
$ws.Range("D6").Value = "'287.26"
$ws.Range("E6").Value = "'  -1.64%  "
$ws.Range("D7").Value = "'0.3661"
$ws.Range("E7").Value = "'  -1.69%  "
$ws.Range("D8").Value = "'48.06"
$ws.Range("E8").Value = "'  -3.76%  "
$ws.Range("D9").Value = "'0.3346"
$ws.Range("E9").Value = "'  -1.55%  "
$ws.Range("D10").Value = "'1.131"
$ws.Range("E10").Value = "'  -1.36%  "
$ws.Range("D11").Value = "'0.07449"
$ws.Range("E11").Value = "'  -1.43%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "'  -0.03%  "
$ws.Range("D13").Value = "'20.90"
$ws.Range("D14").Value = "'6.006"
$ws.Range("E14").Value = "'  -0.55%  "
$ws.Range("D15").Value = "'6.928"
$ws.Range("E15").Value = "'  -0.52%  "
$ws.Range("D16").Value = "'1.573.82"
$ws.Range("E16").Value = "'  +0.42%  "
$ws.Range("E17").Value = "'  -1.01%  "
$ws.Range("D18").Value = "'88.32"
$ws.Range("E18").Value = "'  -2.78%  "
$ws.Range("D19").Value = "'0.06757"
$ws.Range("E19").Value = "'  -0.06%  "
$ws.Range("B20").Value = "'Uniswap"
$ws.Range("C20").Value = "'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'6.433"
$ws.Range("E20").Value = "'  +2.07%  "
$ws.Range("B21").Value = "'Dai"
$ws.Range("C21").Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "'1.002"
$ws.Range("E21").Value = "'  +0.02%  "
$ws.Range("D22").Value = "'16.47"
$ws.Range("E22").Value = "'  +0.69%  "
$ws.Range("D23").Value = "'12.14"
$ws.Range("E23").Value = "'  -0.43%  "
$ws.Range("D24").Value = "'22.487.37"
$ws.Range("E24").Value = "'  +0.00%  "
$ws.Range("D25").Value = "'2.386"
$ws.Range("E25").Value = "'  +0.18%  "
$ws.Range("D26").Value = "'2.634"
$ws.Range("E26").Value = "'  +0.34%  "
$ws.Range("D27").Value = "'152.34"
$ws.Range("E27").Value = "'  +2.11%  "
$ws.Range("D28").Value = "'19.67"
$ws.Range("E28").Value = "'  -2.01%  "
$ws.Range("D29").Value = "'5.001"
$ws.Range("E29").Value = "'  -1.08%  "
$ws.Range("D30").Value = "'124.33"
$ws.Range("E30").Value = "'  -0.79%  "
$ws.Range("D31").Value = "'1.752.77"
$ws.Range("E31").Value = "'  +0.48%  "
$ws.Range("D32").Value = "'1.043"
$ws.Range("E32").Value = "'  -3.97%  "
$ws.Range("D33").Value = "'6.189"
$ws.Range("E33").Value = "'  -0.29%  "
$ws.Range("D34").Value = "'2.001"
$ws.Range("E34").Value = "'  -0.29%  "
$ws.Range("E35").Value = "'  +0.15%  "
$ws.Range("D36").Value = "'0.08288"
$ws.Range("E36").Value = "'  -1.03%  "
$ws.Range("D37").Value = "'0.02443"
$ws.Range("E37").Value = "'  -1.78%  "
$ws.Range("D38").Value = "'0.2267"
$ws.Range("E38").Value = "'  -1.70%  "
$ws.Range("D39").Value = "'0.06487"
$ws.Range("E39").Value = "'  -0.82%  "
$ws.Range("D40").Value = "'5.431"
$ws.Range("E40").Value = "'  -0.41%  "
$ws.Range("D41").Value = "'1.301"
$ws.Range("E41").Value = "'  -2.89%  "
$ws.Range("D42").Value = "'11.31"
$ws.Range("E42").Value = "'  -0.22%  "
$ws.Range("D43").Value = "'0.6306"
$ws.Range("E43").Value = "'  +1.23%  "
$ws.Range("B44").Value = "'EnergySwap"
$ws.Range("C44").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'14.01"
$ws.Range("E44").Value = "'  +0.14%  "
$ws.Range("B45").Value = "'Decentraland"
$ws.Range("C45").Value = "'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.6096"
$ws.Range("E45").Value = "'  +4.12%  "
$ws.Range("B46").Value = "'PancakeSwap"
$ws.Range("C46").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").Value = "'3.760"
$ws.Range("E46").Value = "'  -1.39%  "
$ws.Range("B47").Value = "'NEARProtocol"
$ws.Range("C47").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "'2.057"
$ws.Range("E47").Value = "'  -0.60%  "
$ws.Range("B48").Value = "'Quant"
$ws.Range("C48").Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "'124.81"
$ws.Range("E48").Value = "'  -4.31%  "
$ws.Range("B49").Value = "'EOS"
$ws.Range("C49").Value = "'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D49").Value = "'1.223"
$ws.Range("E49").Value = "'  +0.70%  "
$ws.Range("B50").Value = "'Cronos"
$ws.Range("C50").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.07233"
$ws.Range("E50").Value = "'  -1.31%  "
$ws.Range("B51").Value = "'Aave"
$ws.Range("C51").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'76.87"
$ws.Range("E51").Value = "'  +0.32%  "
